$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-19 Thursday" "2026-02-20 Friday"

Replace-Text "429×7=" "636×2="
Replace-Text "114×3=" "706×8="
Replace-Text "723×9=" "665×4="
Replace-Text "922×8=" "744×3="
Replace-Text "278×3=" "612×3="

Replace-Text "336×6=" "686×6="
Replace-Text "132×9=" "728×3="
Replace-Text "550×8=" "626×3="
Replace-Text "936×5=" "907×4="
Replace-Text "146×8=" "712×5="

Replace-Text "684×3=" "808×4="
Replace-Text "278×8=" "588×6="
Replace-Text "625×2=" "252×8="
Replace-Text "361×7=" "233×6="
Replace-Text "743×7=" "568×6="

Replace-Text "450×3=" "196×2="
Replace-Text "807×2=" "590×5="
Replace-Text "888×4=" "540×5="
Replace-Text "590×3=" "283×8="
Replace-Text "488×2=" "591×7="

Replace-Text "372×5=" "913×8="
Replace-Text "801×5=" "700×8="
Replace-Text "378×8=" "609×4="
Replace-Text "581×3=" "461×3="
Replace-Text "607×4=" "737×4="
